$d = $word.ActiveDocument

# The paragraph "<id>p017r_1</id>" is split across three runs:
#   run 1: "<id>"   (Courier New, color 7f6000, sz 18)
#   run 2: "p017r_1" (plain, color 000000)
#   run 3: "</id>"  (Courier New, color 7f6000, sz 18)
# Merge them into a single run "<id>p017r_1</id>" carrying the formatting
# of the first run, by deleting the text contributed by runs 2 & 3 and
# re-inserting it right after run 1 (inherits run 1's character formatting).

$needle = "<id>p017r_1</id>"
$splitAt = "<id>".Length

$full = $d.Content
$full.Find.ClearFormatting()
$full.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($full.Find.Found) {
    $matchStart = $full.Start
    $matchEnd = $full.End

    $firstRun = $d.Range($matchStart, $matchStart + $splitAt)
    $rest = $d.Range($matchStart + $splitAt, $matchEnd)

    $restText = $rest.Text
    $rest.Delete()
    $firstRun.InsertAfter($restText)
}
